$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 345; this shifts existing rows 345..432 down to 346..433
$ws.Rows.Item(345).Insert()

# Populate the new row 345 with the new data record
$ws.Range("A345").Value = 9
$ws.Range("B345").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C345").Value = "Metropolitana"
$ws.Range("D345").Value = 44722
$ws.Range("E345").Value = 13
$ws.Range("F345").Value = "Fruta"
$ws.Range("G345").Value = 100108
$ws.Range("H345").Value = "Tropicales y subtropicales"
$ws.Range("I345").Value = 100108002
$ws.Range("J345").Value = "Mango"
$ws.Range("K345").Value = "Sin especificar"
$ws.Range("L345").Value = "Primera"
$ws.Range("M345").Value = 500
$ws.Range("N345").Value = 8000
$ws.Range("O345").Value = 9000
$ws.Range("P345").Value = 8400
$ws.Range("Q345").Value = '$/bandeja 4 kilos'
$ws.Range("R345").Value = "Brasil"
$ws.Range("S345").Value = 2100
$ws.Range("T345").Value = 4
